$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 was an empty placeholder row for the "background color" demo.
# Fill it in: label text + a solid red fill (Interior), matching the
# pattern used by the other font/border demo rows above it.
$ws.Range("A21").Value = "background color"
$ws.Range("A21").Interior.Pattern = 1
$ws.Range("A21").Interior.Color = 255

# The "indexed font color" demo (row 23) now uses a theme color
# (Accent1, theme index 4) instead of the legacy palette-indexed color.
$ws.Range("A23").Font.ThemeColor = 5

Write-Host "Applied background-color fill demo edits"
